$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the second half of the table (Ofast, Os, Og, Og funroll loops)
# from columns F:I (rows 1-7) down to a new table in columns A:E (rows 9-15) ---

# Headers for the new table
$ws.Range("A9").Value = $ws.Range("A1").Value2
$ws.Range("B9").Value = $ws.Range("F1").Value2
$ws.Range("C9").Value = $ws.Range("G1").Value2
$ws.Range("D9").Value = $ws.Range("H1").Value2
$ws.Range("E9").Value = $ws.Range("I1").Value2

# Data rows (rows 2-6 -> rows 10-14), column A keeps the run numbers
for ($i = 0; $i -lt 5; $i++) {
    $srcRow = 2 + $i
    $dstRow = 10 + $i
    $ws.Cells.Item($dstRow, 1).Value = $ws.Cells.Item($srcRow, 1).Value2
    $ws.Cells.Item($dstRow, 2).Value = $ws.Cells.Item($srcRow, 6).Value2
    $ws.Cells.Item($dstRow, 3).Value = $ws.Cells.Item($srcRow, 7).Value2
    $ws.Cells.Item($dstRow, 4).Value = $ws.Cells.Item($srcRow, 8).Value2
    $ws.Cells.Item($dstRow, 5).Value = $ws.Cells.Item($srcRow, 9).Value2
}

# Mean row for the new table
$ws.Range("A15").Value = $ws.Range("A7").Value2
$ws.Range("B15").Formula = "=AVERAGE(B10:B14)"
$ws.Range("C15").Formula = "=AVERAGE(C10:C14)"
$ws.Range("D15").Formula = "=AVERAGE(D10:D14)"
$ws.Range("E15").Formula = "=AVERAGE(E10:E14)"

# --- Remove the now-duplicated columns F:I from rows 1-7 ---
$ws.Range("F1:I7").Clear()

# --- Update the sheet's active selection to match the final layout ---
$ws.Range("B9:E15").Select()
